# Auto-generated edit script applying the diff to resum_diari_meteocat sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell used purely to re-copy a known-good ("General", border-only) cell
# style onto percent cells after forcing a text NumberFormat - prevents Excel
# from auto-converting strings like "80%" into a numeric 0.8 with % formatting,
# which would change both the stored type and the style id of the cell.
$styleDonor = $ws.Range("C2")

$ws.Range("E2").Value = "2026-02-06 23:47:40"
$ws.Range("O2").Value = "0.0 °C"
$ws.Range("E3").Value = "2026-02-06 23:47:43"
$ws.Range("E4").Value = "2026-02-06 23:47:45"
$ws.Range("E5").Value = "2026-02-06 23:47:48"
$ws.Range("J5").Value = "998.3 hPa"
$ws.Range("E6").Value = "2026-02-06 23:47:50"
$ws.Range("O6").Value = "15.1 °C"
$ws.Range("E7").Value = "2026-02-06 23:47:53"
$ws.Range("E8").Value = "2026-02-06 23:47:55"
$ws.Range("O8").Value = "9.3 °C"
$ws.Range("E9").Value = "2026-02-06 23:47:58"
$ws.Range("I9").Value = "0.4 mm"
$ws.Range("E10").Value = "2026-02-06 23:48:00"
$ws.Range("E11").Value = "2026-02-06 23:48:02"
$ws.Range("J11").Value = "999.8 hPa"
$ws.Range("O11").Value = "4.8 °C"
$ws.Range("E12").Value = "2026-02-06 23:48:05"
$ws.Range("O12").Value = "12.9 °C"
$ws.Range("E13").Value = "2026-02-06 23:48:07"
$ws.Range("O13").Value = "9.8 °C"
$ws.Range("E14").Value = "2026-02-06 23:48:10"
$ws.Range("O14").Value = "-4.5 °C"
$ws.Range("E15").Value = "2026-02-06 23:48:12"
$ws.Range("J15").Value = "998.5 hPa"
$ws.Range("E16").Value = "2026-02-06 23:48:15"
$ws.Range("E17").Value = "2026-02-06 23:48:17"
$ws.Range("I17").Value = "0.7 mm"
$ws.Range("J17").Value = "999.8 hPa"
$ws.Range("O17").Value = "6.1 °C"
$ws.Range("E18").Value = "2026-02-06 23:48:20"
$ws.Range("E19").Value = "2026-02-06 23:48:22"
$r = $ws.Range("H19")
$r.NumberFormat = "@"
$r.Value = "80%"
$styleDonor.Copy()
$r.PasteSpecial(-4122)

$ws.Range("J19").Value = "1000.7 hPa"
$ws.Range("E20").Value = "2026-02-06 23:48:25"
$ws.Range("O20").Value = "-2.3 °C"
$ws.Range("E21").Value = "2026-02-06 23:48:27"
$ws.Range("E22").Value = "2026-02-06 23:48:30"
$ws.Range("O22").Value = "9.8 °C"
$ws.Range("E23").Value = "2026-02-06 23:48:32"
$r = $ws.Range("H23")
$r.NumberFormat = "@"
$r.Value = "86%"
$styleDonor.Copy()
$r.PasteSpecial(-4122)

$ws.Range("E24").Value = "2026-02-06 23:48:35"
$r = $ws.Range("H24")
$r.NumberFormat = "@"
$r.Value = "69%"
$styleDonor.Copy()
$r.PasteSpecial(-4122)

$ws.Range("O24").Value = "12.6 °C"
$ws.Range("E25").Value = "2026-02-06 23:48:37"
$ws.Range("J25").Value = "999.4 hPa"
$ws.Range("O25").Value = "4.2 °C"
$ws.Range("E26").Value = "2026-02-06 23:48:40"
$ws.Range("E27").Value = "2026-02-06 23:48:42"
$ws.Range("J27").Value = "998.4 hPa"
$ws.Range("E28").Value = "2026-02-06 23:48:44"
$ws.Range("E29").Value = "2026-02-06 23:48:47"
$ws.Range("O29").Value = "12.1 °C"
$ws.Range("E30").Value = "2026-02-06 23:48:49"
$ws.Range("E31").Value = "2026-02-06 23:48:51"
$ws.Range("I31").Value = "5.2 mm"
$ws.Range("J31").Value = "1000.2 hPa"
$ws.Range("O31").Value = "7.0 °C"
$ws.Range("E32").Value = "2026-02-06 23:48:53"
$ws.Range("J32").Value = "999.6 hPa"
$ws.Range("O32").Value = "15.0 °C"
$ws.Range("E33").Value = "2026-02-06 23:48:56"
$ws.Range("E34").Value = "2026-02-06 23:48:58"
$ws.Range("E35").Value = "2026-02-06 23:49:00"
$ws.Range("N35").Value = "-4.1 °C 23:20 TU"
$ws.Range("E36").Value = "2026-02-06 23:49:03"
$ws.Range("J36").Value = "1000.8 hPa"
$ws.Range("N36").Value = "5.9 °C 23:20 TU"
$ws.Range("O36").Value = "12.0 °C"
